$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C (Förändrad) holds a date serial value that was updated from
# 45179 (2023-09-10) to 45180 (2023-09-11) for every data row (rows 2-97).
$ws.Range("C2:C97").Value = 45180
